$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric stat corrections (GP/W/L/... and rank columns)
$ws.Cells.Item(2, 30).Value = 19
$ws.Cells.Item(2, 32).Value = 12
$ws.Cells.Item(2, 33).Value = 12
$ws.Cells.Item(2, 34).Value = 9
$ws.Cells.Item(2, 36).Value = 20
$ws.Cells.Item(2, 40).Value = 9
$ws.Cells.Item(2, 42).Value = 22
$ws.Cells.Item(2, 48).Value = 20
$ws.Cells.Item(2, 51).Value = 13
$ws.Cells.Item(2, 53).Value = 23
$ws.Cells.Item(2, 55).Value = 14
$ws.Cells.Item(3, 4).Value = 41
$ws.Cells.Item(3, 6).Value = 27
$ws.Cells.Item(3, 7).Value = 0.341
$ws.Cells.Item(3, 9).Value = 36.5
$ws.Cells.Item(3, 10).Value = 82.8
$ws.Cells.Item(3, 11).Value = 0.441
$ws.Cells.Item(3, 12).Value = 6.3
$ws.Cells.Item(3, 13).Value = 18.8
$ws.Cells.Item(3, 14).Value = 0.334
$ws.Cells.Item(3, 15).Value = 16
$ws.Cells.Item(3, 17).Value = 0.763
$ws.Cells.Item(3, 18).Value = 11.2
$ws.Cells.Item(3, 21).Value = 19.8
$ws.Cells.Item(3, 22).Value = 15.5
$ws.Cells.Item(3, 26).Value = 21.9
$ws.Cells.Item(3, 28).Value = 95.3
$ws.Cells.Item(3, 30).Value = 7
$ws.Cells.Item(3, 32).Value = 26
$ws.Cells.Item(3, 35).Value = 22
$ws.Cells.Item(3, 36).Value = 18
$ws.Cells.Item(3, 37).Value = 20
$ws.Cells.Item(3, 41).Value = 24
$ws.Cells.Item(3, 43).Value = 13
$ws.Cells.Item(3, 46).Value = 20
$ws.Cells.Item(3, 48).Value = 22
$ws.Cells.Item(3, 50).Value = 15
$ws.Cells.Item(3, 53).Value = 28
$ws.Cells.Item(4, 30).Value = 28
$ws.Cells.Item(4, 42).Value = 6
$ws.Cells.Item(4, 43).Value = 15
$ws.Cells.Item(4, 49).Value = 21
$ws.Cells.Item(5, 32).Value = 22
$ws.Cells.Item(5, 44).Value = 26
$ws.Cells.Item(5, 47).Value = 26
$ws.Cells.Item(6, 30).Value = 19
$ws.Cells.Item(6, 32).Value = 15
$ws.Cells.Item(6, 43).Value = 10
$ws.Cells.Item(6, 45).Value = 11
$ws.Cells.Item(6, 47).Value = 13
$ws.Cells.Item(7, 30).Value = 8
$ws.Cells.Item(7, 31).Value = 21
$ws.Cells.Item(7, 32).Value = 22
$ws.Cells.Item(7, 33).Value = 22
$ws.Cells.Item(7, 38).Value = 20
$ws.Cells.Item(7, 40).Value = 10
$ws.Cells.Item(7, 48).Value = 14
$ws.Cells.Item(7, 49).Value = 22
$ws.Cells.Item(7, 55).Value = 25
$ws.Cells.Item(8, 33).Value = 9
$ws.Cells.Item(8, 39).Value = 9
$ws.Cells.Item(8, 41).Value = 23
$ws.Cells.Item(9, 4).Value = 39
$ws.Cells.Item(9, 6).Value = 19
$ws.Cells.Item(9, 7).Value = 0.513
$ws.Cells.Item(9, 13).Value = 22.6
$ws.Cells.Item(9, 14).Value = 0.361
$ws.Cells.Item(9, 15).Value = 18.7
$ws.Cells.Item(9, 16).Value = 25.6
$ws.Cells.Item(9, 20).Value = 46
$ws.Cells.Item(9, 21).Value = 22.2
$ws.Cells.Item(9, 22).Value = 14.9
$ws.Cells.Item(9, 24).Value = 6
$ws.Cells.Item(9, 29).Value = 0.9
$ws.Cells.Item(9, 30).Value = 19
$ws.Cells.Item(9, 32).Value = 12
$ws.Cells.Item(9, 33).Value = 12
$ws.Cells.Item(9, 36).Value = 7
$ws.Cells.Item(9, 38).Value = 9
$ws.Cells.Item(9, 39).Value = 10
$ws.Cells.Item(9, 40).Value = 14
$ws.Cells.Item(9, 46).Value = 4
$ws.Cells.Item(9, 47).Value = 12
$ws.Cells.Item(9, 48).Value = 13
$ws.Cells.Item(9, 49).Value = 20
$ws.Cells.Item(9, 50).Value = 4
$ws.Cells.Item(9, 52).Value = 27
$ws.Cells.Item(9, 55).Value = 13
$ws.Cells.Item(10, 30).Value = 8
$ws.Cells.Item(10, 48).Value = 23
$ws.Cells.Item(10, 52).Value = 15
$ws.Cells.Item(11, 46).Value = 5
$ws.Cells.Item(11, 54).Value = 9
$ws.Cells.Item(12, 40).Value = 23
$ws.Cells.Item(12, 45).Value = 4
$ws.Cells.Item(12, 52).Value = 14
$ws.Cells.Item(13, 30).Value = 19
$ws.Cells.Item(13, 35).Value = 23
$ws.Cells.Item(13, 40).Value = 13
$ws.Cells.Item(13, 48).Value = 21
$ws.Cells.Item(13, 50).Value = 5
$ws.Cells.Item(13, 51).Value = 13
$ws.Cells.Item(14, 38).Value = 11
$ws.Cells.Item(14, 45).Value = 12
$ws.Cells.Item(14, 46).Value = 15
$ws.Cells.Item(14, 50).Value = 14
$ws.Cells.Item(15, 4).Value = 40
$ws.Cells.Item(15, 5).Value = 15
$ws.Cells.Item(15, 7).Value = 0.375
$ws.Cells.Item(15, 10).Value = 83.8
$ws.Cells.Item(15, 11).Value = 0.44
$ws.Cells.Item(15, 13).Value = 25.1
$ws.Cells.Item(15, 14).Value = 0.365
$ws.Cells.Item(15, 15).Value = 17.4
$ws.Cells.Item(15, 16).Value = 23.1
$ws.Cells.Item(15, 17).Value = 0.751
$ws.Cells.Item(15, 18).Value = 10.1
$ws.Cells.Item(15, 19).Value = 33
$ws.Cells.Item(15, 20).Value = 43.1
$ws.Cells.Item(15, 22).Value = 15.7
$ws.Cells.Item(15, 24).Value = 5.8
$ws.Cells.Item(15, 26).Value = 20.5
$ws.Cells.Item(15, 28).Value = 100.3
$ws.Cells.Item(15, 29).Value = -5.5
$ws.Cells.Item(15, 31).Value = 21
$ws.Cells.Item(15, 32).Value = 22
$ws.Cells.Item(15, 33).Value = 22
$ws.Cells.Item(15, 37).Value = 22
$ws.Cells.Item(15, 40).Value = 11
$ws.Cells.Item(15, 42).Value = 14
$ws.Cells.Item(15, 43).Value = 18
$ws.Cells.Item(15, 45).Value = 10
$ws.Cells.Item(15, 46).Value = 17
$ws.Cells.Item(15, 47).Value = 10
$ws.Cells.Item(15, 48).Value = 24
$ws.Cells.Item(15, 52).Value = 16
$ws.Cells.Item(15, 53).Value = 25
$ws.Cells.Item(15, 55).Value = 26
$ws.Cells.Item(16, 30).Value = 19
$ws.Cells.Item(16, 32).Value = 12
$ws.Cells.Item(16, 33).Value = 12
$ws.Cells.Item(16, 34).Value = 15
$ws.Cells.Item(16, 35).Value = 16
$ws.Cells.Item(16, 37).Value = 12
$ws.Cells.Item(16, 41).Value = 25
$ws.Cells.Item(16, 42).Value = 24
$ws.Cells.Item(16, 43).Value = 21
$ws.Cells.Item(16, 46).Value = 16
$ws.Cells.Item(16, 53).Value = 23
$ws.Cells.Item(17, 30).Value = 8
$ws.Cells.Item(18, 4).Value = 39
$ws.Cells.Item(18, 6).Value = 32
$ws.Cells.Item(18, 7).Value = 0.179
$ws.Cells.Item(18, 10).Value = 82.40000000000001
$ws.Cells.Item(18, 11).Value = 0.421
$ws.Cells.Item(18, 13).Value = 20.6
$ws.Cells.Item(18, 14).Value = 0.35
$ws.Cells.Item(18, 15).Value = 14.9
$ws.Cells.Item(18, 17).Value = 0.762
$ws.Cells.Item(18, 20).Value = 41.4
$ws.Cells.Item(18, 22).Value = 15.9
$ws.Cells.Item(18, 23).Value = 7
$ws.Cells.Item(18, 25).Value = 5.3
$ws.Cells.Item(18, 27).Value = 20
$ws.Cells.Item(18, 28).Value = 91.59999999999999
$ws.Cells.Item(18, 29).Value = -8.699999999999999
$ws.Cells.Item(18, 30).Value = 19
$ws.Cells.Item(18, 34).Value = 2
$ws.Cells.Item(18, 38).Value = 18
$ws.Cells.Item(18, 43).Value = 14
$ws.Cells.Item(18, 46).Value = 24
$ws.Cells.Item(18, 48).Value = 26
$ws.Cells.Item(18, 49).Value = 25
$ws.Cells.Item(18, 55).Value = 29
$ws.Cells.Item(19, 30).Value = 8
$ws.Cells.Item(19, 34).Value = 22
$ws.Cells.Item(19, 52).Value = 1
$ws.Cells.Item(20, 30).Value = 19
$ws.Cells.Item(20, 31).Value = 21
$ws.Cells.Item(20, 33).Value = 21
$ws.Cells.Item(20, 43).Value = 12
$ws.Cells.Item(21, 30).Value = 8
$ws.Cells.Item(21, 31).Value = 21
$ws.Cells.Item(21, 32).Value = 22
$ws.Cells.Item(21, 33).Value = 22
$ws.Cells.Item(21, 42).Value = 29
$ws.Cells.Item(21, 43).Value = 11
$ws.Cells.Item(22, 4).Value = 40
$ws.Cells.Item(22, 5).Value = 30
$ws.Cells.Item(22, 7).Value = 0.75
$ws.Cells.Item(22, 8).Value = 48.3
$ws.Cells.Item(22, 10).Value = 82.8
$ws.Cells.Item(22, 14).Value = 0.341
$ws.Cells.Item(22, 16).Value = 25.9
$ws.Cells.Item(22, 17).Value = 0.8149999999999999
$ws.Cells.Item(22, 18).Value = 11.1
$ws.Cells.Item(22, 19).Value = 35.7
$ws.Cells.Item(22, 23).Value = 8.1
$ws.Cells.Item(22, 24).Value = 6.2
$ws.Cells.Item(22, 25).Value = 3.9
$ws.Cells.Item(22, 29).Value = 7.1
$ws.Cells.Item(22, 31).Value = 4
$ws.Cells.Item(22, 34).Value = 22
$ws.Cells.Item(22, 36).Value = 19
$ws.Cells.Item(22, 40).Value = 25
$ws.Cells.Item(22, 48).Value = 25
$ws.Cells.Item(22, 52).Value = 26
$ws.Cells.Item(23, 4).Value = 40
$ws.Cells.Item(23, 5).Value = 10
$ws.Cells.Item(23, 7).Value = 0.25
$ws.Cells.Item(23, 9).Value = 36.5
$ws.Cells.Item(23, 10).Value = 82.90000000000001
$ws.Cells.Item(23, 11).Value = 0.441
$ws.Cells.Item(23, 12).Value = 7.3
$ws.Cells.Item(23, 13).Value = 21.5
$ws.Cells.Item(23, 14).Value = 0.341
$ws.Cells.Item(23, 15).Value = 15.8
$ws.Cells.Item(23, 16).Value = 21
$ws.Cells.Item(23, 17).Value = 0.751
$ws.Cells.Item(23, 19).Value = 33.2
$ws.Cells.Item(23, 20).Value = 42.6
$ws.Cells.Item(23, 25).Value = 6
$ws.Cells.Item(23, 26).Value = 20.3
$ws.Cells.Item(23, 27).Value = 18.8
$ws.Cells.Item(23, 28).Value = 96.09999999999999
$ws.Cells.Item(23, 29).Value = -6.1
$ws.Cells.Item(23, 34).Value = 4
$ws.Cells.Item(23, 35).Value = 21
$ws.Cells.Item(23, 36).Value = 17
$ws.Cells.Item(23, 40).Value = 24
$ws.Cells.Item(23, 41).Value = 26
$ws.Cells.Item(23, 42).Value = 25
$ws.Cells.Item(23, 44).Value = 27
$ws.Cells.Item(23, 46).Value = 21
$ws.Cells.Item(23, 48).Value = 19
$ws.Cells.Item(23, 53).Value = 29
$ws.Cells.Item(24, 30).Value = 8
$ws.Cells.Item(24, 34).Value = 5
$ws.Cells.Item(24, 38).Value = 21
$ws.Cells.Item(24, 47).Value = 11
$ws.Cells.Item(24, 55).Value = 30
$ws.Cells.Item(25, 4).Value = 39
$ws.Cells.Item(25, 5).Value = 22
$ws.Cells.Item(25, 7).Value = 0.5639999999999999
$ws.Cells.Item(25, 9).Value = 38.4
$ws.Cells.Item(25, 10).Value = 84.90000000000001
$ws.Cells.Item(25, 11).Value = 0.452
$ws.Cells.Item(25, 12).Value = 9.5
$ws.Cells.Item(25, 13).Value = 26.2
$ws.Cells.Item(25, 14).Value = 0.361
$ws.Cells.Item(25, 16).Value = 23.1
$ws.Cells.Item(25, 17).Value = 0.751
$ws.Cells.Item(25, 18).Value = 11.8
$ws.Cells.Item(25, 20).Value = 43.7
$ws.Cells.Item(25, 21).Value = 19
$ws.Cells.Item(25, 22).Value = 15.1
$ws.Cells.Item(25, 23).Value = 8.6
$ws.Cells.Item(25, 25).Value = 4.3
$ws.Cells.Item(25, 28).Value = 103.6
$ws.Cells.Item(25, 29).Value = 2.4
$ws.Cells.Item(25, 30).Value = 19
$ws.Cells.Item(25, 33).Value = 10
$ws.Cells.Item(25, 36).Value = 8
$ws.Cells.Item(25, 37).Value = 13
$ws.Cells.Item(25, 40).Value = 15
$ws.Cells.Item(25, 41).Value = 15
$ws.Cells.Item(25, 42).Value = 13
$ws.Cells.Item(25, 43).Value = 20
$ws.Cells.Item(25, 44).Value = 10
$ws.Cells.Item(25, 48).Value = 16
$ws.Cells.Item(25, 54).Value = 10
$ws.Cells.Item(25, 55).Value = 11
$ws.Cells.Item(26, 30).Value = 8
$ws.Cells.Item(26, 31).Value = 2
$ws.Cells.Item(26, 33).Value = 2
$ws.Cells.Item(26, 45).Value = 5
$ws.Cells.Item(26, 50).Value = 16
$ws.Cells.Item(27, 4).Value = 38
$ws.Cells.Item(27, 6).Value = 24
$ws.Cells.Item(27, 7).Value = 0.368
$ws.Cells.Item(27, 9).Value = 37.6
$ws.Cells.Item(27, 10).Value = 83.59999999999999
$ws.Cells.Item(27, 11).Value = 0.45
$ws.Cells.Item(27, 14).Value = 0.349
$ws.Cells.Item(27, 15).Value = 19.6
$ws.Cells.Item(27, 18).Value = 11.6
$ws.Cells.Item(27, 21).Value = 20.1
$ws.Cells.Item(27, 22).Value = 14.9
$ws.Cells.Item(27, 23).Value = 7.8
$ws.Cells.Item(27, 27).Value = 22.5
$ws.Cells.Item(27, 28).Value = 101.7
$ws.Cells.Item(27, 29).Value = -2.1
$ws.Cells.Item(27, 30).Value = 28
$ws.Cells.Item(27, 32).Value = 20
$ws.Cells.Item(27, 35).Value = 15
$ws.Cells.Item(27, 37).Value = 15
$ws.Cells.Item(27, 43).Value = 8
$ws.Cells.Item(27, 47).Value = 24
$ws.Cells.Item(27, 48).Value = 12
$ws.Cells.Item(27, 52).Value = 28
$ws.Cells.Item(27, 53).Value = 5
$ws.Cells.Item(28, 4).Value = 40
$ws.Cells.Item(28, 5).Value = 31
$ws.Cells.Item(28, 7).Value = 0.775
$ws.Cells.Item(28, 9).Value = 40.8
$ws.Cells.Item(28, 10).Value = 83
$ws.Cells.Item(28, 11).Value = 0.492
$ws.Cells.Item(28, 13).Value = 20.5
$ws.Cells.Item(28, 14).Value = 0.399
$ws.Cells.Item(28, 15).Value = 14.9
$ws.Cells.Item(28, 17).Value = 0.774
$ws.Cells.Item(28, 18).Value = 9.1
$ws.Cells.Item(28, 19).Value = 33.8
$ws.Cells.Item(28, 20).Value = 42.9
$ws.Cells.Item(28, 24).Value = 4.7
$ws.Cells.Item(28, 25).Value = 4.6
$ws.Cells.Item(28, 28).Value = 104.6
$ws.Cells.Item(28, 29).Value = 7.7
$ws.Cells.Item(28, 31).Value = 2
$ws.Cells.Item(28, 36).Value = 16
$ws.Cells.Item(28, 38).Value = 10
$ws.Cells.Item(28, 42).Value = 30
$ws.Cells.Item(28, 43).Value = 7
$ws.Cells.Item(28, 45).Value = 6
$ws.Cells.Item(28, 46).Value = 19
$ws.Cells.Item(28, 50).Value = 16
$ws.Cells.Item(28, 51).Value = 12
$ws.Cells.Item(28, 52).Value = 2
$ws.Cells.Item(29, 4).Value = 38
$ws.Cells.Item(29, 6).Value = 18
$ws.Cells.Item(29, 7).Value = 0.526
$ws.Cells.Item(29, 8).Value = 48.7
$ws.Cells.Item(29, 9).Value = 35.7
$ws.Cells.Item(29, 10).Value = 82
$ws.Cells.Item(29, 11).Value = 0.435
$ws.Cells.Item(29, 13).Value = 22.2
$ws.Cells.Item(29, 14).Value = 0.363
$ws.Cells.Item(29, 15).Value = 19.3
$ws.Cells.Item(29, 16).Value = 25
$ws.Cells.Item(29, 17).Value = 0.771
$ws.Cells.Item(29, 19).Value = 31.1
$ws.Cells.Item(29, 20).Value = 43
$ws.Cells.Item(29, 21).Value = 19.9
$ws.Cells.Item(29, 23).Value = 7.1
$ws.Cells.Item(29, 25).Value = 5
$ws.Cells.Item(29, 27).Value = 22.6
$ws.Cells.Item(29, 28).Value = 98.7
$ws.Cells.Item(29, 29).Value = 2.7
$ws.Cells.Item(29, 30).Value = 28
$ws.Cells.Item(29, 32).Value = 10
$ws.Cells.Item(29, 34).Value = 7
$ws.Cells.Item(29, 38).Value = 12
$ws.Cells.Item(29, 40).Value = 12
$ws.Cells.Item(29, 43).Value = 9
$ws.Cells.Item(29, 44).Value = 9
$ws.Cells.Item(29, 46).Value = 18
$ws.Cells.Item(29, 47).Value = 27
$ws.Cells.Item(29, 49).Value = 24
$ws.Cells.Item(29, 53).Value = 4
$ws.Cells.Item(29, 55).Value = 10
$ws.Cells.Item(30, 32).Value = 28
$ws.Cells.Item(30, 33).Value = 27
$ws.Cells.Item(30, 37).Value = 23
$ws.Cells.Item(30, 41).Value = 22
$ws.Cells.Item(30, 42).Value = 21
$ws.Cells.Item(30, 43).Value = 17
$ws.Cells.Item(30, 46).Value = 25
$ws.Cells.Item(31, 30).Value = 19
$ws.Cells.Item(31, 32).Value = 15
$ws.Cells.Item(31, 35).Value = 17
$ws.Cells.Item(31, 37).Value = 16
$ws.Cells.Item(31, 42).Value = 23

# Date column (BF): format "M-D-YYYY-YY" -> ISO "YYYY-MM-DD", kept as text
$ws.Cells.Item(2, 58).Value = "'2014-01-19"
$ws.Cells.Item(3, 58).Value = "'2014-01-19"
$ws.Cells.Item(4, 58).Value = "'2014-01-19"
$ws.Cells.Item(5, 58).Value = "'2014-01-19"
$ws.Cells.Item(6, 58).Value = "'2014-01-19"
$ws.Cells.Item(7, 58).Value = "'2014-01-19"
$ws.Cells.Item(8, 58).Value = "'2014-01-19"
$ws.Cells.Item(9, 58).Value = "'2014-01-19"
$ws.Cells.Item(10, 58).Value = "'2014-01-19"
$ws.Cells.Item(11, 58).Value = "'2014-01-19"
$ws.Cells.Item(12, 58).Value = "'2014-01-19"
$ws.Cells.Item(13, 58).Value = "'2014-01-19"
$ws.Cells.Item(14, 58).Value = "'2014-01-19"
$ws.Cells.Item(15, 58).Value = "'2014-01-19"
$ws.Cells.Item(16, 58).Value = "'2014-01-19"
$ws.Cells.Item(17, 58).Value = "'2014-01-19"
$ws.Cells.Item(18, 58).Value = "'2014-01-19"
$ws.Cells.Item(19, 58).Value = "'2014-01-19"
$ws.Cells.Item(20, 58).Value = "'2014-01-19"
$ws.Cells.Item(21, 58).Value = "'2014-01-19"
$ws.Cells.Item(22, 58).Value = "'2014-01-19"
$ws.Cells.Item(23, 58).Value = "'2014-01-19"
$ws.Cells.Item(24, 58).Value = "'2014-01-19"
$ws.Cells.Item(25, 58).Value = "'2014-01-19"
$ws.Cells.Item(26, 58).Value = "'2014-01-19"
$ws.Cells.Item(27, 58).Value = "'2014-01-19"
$ws.Cells.Item(28, 58).Value = "'2014-01-19"
$ws.Cells.Item(29, 58).Value = "'2014-01-19"
$ws.Cells.Item(30, 58).Value = "'2014-01-19"
$ws.Cells.Item(31, 58).Value = "'2014-01-19"

# Re-apply default style so the text-forcing quote-prefix does not linger on the cell
$ws.Range("BF2:BF31").Style = "Normal"